$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Domicilio"
$ws.Range("B5").Value = "Calle"

$ws.Range("A6").Value = "Ciudad"
$ws.Range("B6").Value = "Monto"

$ws.Range("A7").Value = "Programa académico"
$ws.Range("B7").Value = "Escuela"
